# Revise antibody template 1.2
#
# Applies the "Revise antibody template 1.2" edit:
#  - Instructions sheet: rewritten help text, new columns A/B layout
#  - Antibodies sheet: new columns (Light chain, Heavy chain germline,
#    Structural data), reworked sample data, two extra data validations
#  - Terminology sheet: new columns (Light chain, Heavy chain germline)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: apply the bold "header" style used throughout this workbook by
# copying the format from a cell that already carries it (style index 1
# in the original file), so we don't fork a brand-new font/style entry.
# ---------------------------------------------------------------------
function Set-HeaderStyle($ws, $srcA1, [string[]]$destA1List) {
    $ws.Range($srcA1).Copy()
    foreach ($dest in $destA1List) {
        $ws.Range($dest).PasteSpecial(-4122)  # xlPasteFormats
    }
}

# =======================================================================
# Sheet: Instructions
# =======================================================================
$wsI = $wb.Worksheets.Item("Instructions")
$wsI.Unprotect()

# Column widths
$wsI.Columns.Item(1).ColumnWidth = 17.1666666667   # -> width 18
$wsI.Columns.Item(2).ColumnWidth = 69.1666666667   # -> width 70

# Wipe out the old rows 2:10 ("(blank)", intro text, "Columns:", bullet
# list) and re-insert 21 fresh blank rows (2:22) so the new, longer
# content can be written without old leftovers shifting underneath it.
$wsI.Rows("2:10").Delete()
$wsI.Rows("2:22").Insert()
# Insert() inherits the row-above's formatting (bold, from row 1) -
# strip it back to the default so only the headings we re-bold below
# end up styled.
$wsI.Range("A2:B22").ClearFormats()

# Row 1
$wsI.Range("B1").Value = ""

# Row 2
$wsI.Range("A2").Value = "Version 1.2"

# Row 3 (blank)

# Row 4
$wsI.Range("A4").Value = "Add your antibodies to the 'Antibodies' sheet. Do not edit the other sheets."

# Row 6
$wsI.Range("A6").Value = "Antibody name"
$wsI.Range("B6").Value = "Your preferred code name for the antibody"

# Row 7
$wsI.Range("A7").Value = "Host"
$wsI.Range("B7").Value = "Specify the host species that is the source of the antibody"

# Row 8
$wsI.Range("A8").Value = "Isotype"
$wsI.Range("B8").Value = "Specify the antibody isotype, if known"

# Row 9
$wsI.Range("A9").Value = "Light chain"
$wsI.Range("B9").Value = "Specify the antibody light chain, if known (kappa or lambda)"

# Row 10
$wsI.Range("A10").Value = "Heavy chain germline"
$wsI.Range("B10").Value = "Specify the antibody heavy chain germline gene, if known"

# Row 11
$wsI.Range("A11").Value = "Antibody details"
$wsI.Range("B11").Value = "Measurements or characteristics of the antibody."

# Row 12
$wsI.Range("B12").Value = "This column is optional, and meant to capture data you might have on the antibody."

# Row 13
$wsI.Range("B13").Value = "These data will not be released to the partner reference labs that will perform the analyses."

# Row 14
$wsI.Range("B14").Value = "For example:"

# Row 15
$wsI.Range("B15").Value = "- Affinity: Spike protein binding affinity; inhibition of ACE2 binding; ELISA for Spike "

# Row 16
$wsI.Range("B16").Value = "- Neutralization: IC50 value"

# Row 17
$wsI.Range("B17").Value = "- Neutralization assay platform"

# Row 18
$wsI.Range("B18").Value = "- Epitope: Binning or competition data"

# Row 19
$wsI.Range("A19").Value = "Structural data"
$wsI.Range("B19").Value = "Would you like structural analyses of this antibody?"

# Row 20
$wsI.Range("B20").Value = "If no, leave blank."

# Row 21
$wsI.Range("B21").Value = "If yes, rank the antibodies in order of priority, starting with '1' for the highest priority."

# Row 22
$wsI.Range("A22").Value = "Antibody comment"
$wsI.Range("B22").Value = "Please provide any other details about the antibody."

# Bold "header" style on column A label cells (A1 already has it)
Set-HeaderStyle $wsI "A1" @("A6","A7","A8","A9","A10","A11","A19","A22")

$wsI.Protect()

Write-Host "Instructions sheet updated"

# =======================================================================
# Sheet: Antibodies
# =======================================================================
$wsA = $wb.Worksheets.Item("Antibodies")

# Insert two new columns at D:E (Light chain, Heavy chain germline) -
# this pushes the old D (Antibody details) to F and old E (Antibody
# comment) to H - then insert one more new column at G (Structural data)
# between them.
$wsA.Columns("D:E").Insert()
$wsA.Columns("G").Insert()

# Column widths: A,B,C=15 (unchanged), D=15, E=20, F=16 (was D, unchanged),
# G=15, H=16 (was E, unchanged)
$wsA.Columns.Item(4).ColumnWidth = 14.1666666667   # -> width 15
$wsA.Columns.Item(5).ColumnWidth = 19.1666666667   # -> width 20
$wsA.Columns.Item(7).ColumnWidth = 14.1666666667   # -> width 15

# Header row (D1/E1/G1 already inherited the bold header style from the
# column insert; just fill in the text)
$wsA.Range("D1").Value = "Light chain"
$wsA.Range("E1").Value = "Heavy chain germline"
$wsA.Range("G1").Value = "Structural data"

# --- Data rows -----------------------------------------------------
# Row 2
$wsA.Range("A2").Value = "VD-Crotty 1"
$wsA.Range("D2").Value = "kappa"
$wsA.Range("E2").Value = "IGHV1-8"
$wsA.Range("F2").Value = "Spike protein binding affinity"

# Row 3
$wsA.Range("A3").Value = "VD-Crotty 1"
$wsA.Range("B3").Value = "Homo sapiens"
$wsA.Range("D3").Value = "lambda"
$wsA.Range("E3").Value = "IGHV1-18"
$wsA.Range("F3").Value = "inhibition of ACE2 binding"
$wsA.Range("G3").Value = "3"

# Row 4
$wsA.Range("C4").Value = "IgG"
$wsA.Range("E4").Value = "IGHV2-5"
$wsA.Range("F4").Value = "ELISA for Spike"
$wsA.Range("G4").Value = "6"

# Row 5
$wsA.Range("A5").Value = "VD-Crotty 4"
$wsA.Range("B5").Value = ""
$wsA.Range("C5").Value = "IgG2a"
$wsA.Range("E5").Value = "IGHV3-7"
$wsA.Range("F5").Value = "IC50 value"

# Row 6
$wsA.Range("A6").Value = "VD-Crotty 5"
$wsA.Range("B6").Value = "Mus musculus"
$wsA.Range("C6").Value = "IggA1"
$wsA.Range("D6").Value = "kappa"
$wsA.Range("E6").Value = "IGHV3-11"
$wsA.Range("F6").Value = "Neutralization assay platform"
$wsA.Range("G6").Value = "5"
$wsA.Range("H6").Value = "A comment"

# Row 7
$wsA.Range("A7").Value = "VD-Crotty 6"
$wsA.Range("B7").Value = "Mus musculus"
$wsA.Range("C7").Value = "IgA"
$wsA.Range("D7").Value = "kapa"
$wsA.Range("F7").Value = "Epitope binning data"

# Row 8
$wsA.Range("A8").Value = "VD-Crotty 7"
$wsA.Range("E8").Value = "IGVH1-8"
$wsA.Range("F8").Value = "Epitope competition data"
$wsA.Range("G8").Value = "1"

# Row 9
$wsA.Range("A9").Value = "VD-Crotty 8"
$wsA.Range("B9").Value = "Mus musculus"
$wsA.Range("C9").Value = "IgA2"
$wsA.Range("D9").Value = "lambda"
$wsA.Range("G9").Value = "top"
$wsA.Range("H9").Value = "Another comment"

# Row 10
$wsA.Range("A10").Value = "VD-Crotty 9"
$wsA.Range("C10").Value = "IgG1"

# Row 11 (new)
$wsA.Range("A11").Value = "VD-Crotty 10"
$wsA.Range("B11").Value = "Mus musculus"
$wsA.Range("C11").Value = "IgM"
$wsA.Range("G11").Value = "2"

# --- Data validations ------------------------------------------------
$wsA.Range("D2:D100").Validation.Add(3, 1, 1, "Terminology!`$C`$2:`$C`$3")
$wsA.Range("E2:E100").Validation.Add(3, 1, 1, "Terminology!`$D`$2:`$D`$12")

Write-Host "Antibodies sheet updated"

# =======================================================================
# Sheet: Terminology
# =======================================================================
$wsT = $wb.Worksheets.Item("Terminology")
$wsT.Unprotect()

# New columns C (Light chain) and D (Heavy chain germline), appended
# after the existing A/B columns.
$wsT.Columns.Item(3).ColumnWidth = 14.1666666667   # -> width 15
$wsT.Columns.Item(4).ColumnWidth = 19.1666666667   # -> width 20

# Header row
$wsT.Range("C1").Value = "Light chain"
$wsT.Range("D1").Value = "Heavy chain germline"
Set-HeaderStyle $wsT "A1" @("C1","D1")

# Data rows
$wsT.Range("C2").Value = "kappa"
$wsT.Range("D2").Value = "IGHV1-8"

$wsT.Range("C3").Value = "lambda"
$wsT.Range("D3").Value = "IGHV1-18"

$wsT.Range("D4").Value = "IGHV2-5"
$wsT.Range("D5").Value = "IGHV3-7"
$wsT.Range("D6").Value = "IGHV3-11"
$wsT.Range("D7").Value = "IGHV3-21"
$wsT.Range("D8").Value = "IGHV3-23"
$wsT.Range("D9").Value = "IGHV4-39"
$wsT.Range("D10").Value = "IGHV4-59"
$wsT.Range("D11").Value = "IGHV5-51"
$wsT.Range("D12").Value = "IGHV6-1"

# Rows 4-15 of column C, and rows 13-15 of column D, exist as blank
# cells in the target - touch them so the used range / cell nodes line
# up.
$wsT.Range("C4:C15").Value = ""
$wsT.Range("D13:D15").Value = ""

$wsT.Protect()

Write-Host "Terminology sheet updated"
